# Extend each of the 4 worksheets (svr, gbr, rf, lr) from columns B:F (5 models)
# to columns B:K (10 models), updating the existing B:F values and filling in the
# new G:K values, per the updated experiment results.

$wb = $excel.ActiveWorkbook

### Sheet 1 ###
$ws = $wb.Worksheets.Item(1)

# Header row (row 1): copy formatting from F1 (bordered/bold/centered style)
# into the new header cells G1:K1, then set their values 6..10.
$ws.Range("F1").Copy()
$ws.Range("G1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10

# Row 2 (test_MAPE): updated values for all 10 models.
$ws.Range("B2").Value = 28.98029672317424
$ws.Range("C2").Value = 36.70850919668717
$ws.Range("D2").Value = 37.96986392644682
$ws.Range("E2").Value = 31.22592845419224
$ws.Range("F2").Value = 43.61157316185502
$ws.Range("G2").Value = 30.27784069838451
$ws.Range("H2").Value = 31.30217745473083
$ws.Range("I2").Value = 39.61213871867896
$ws.Range("J2").Value = 42.51503485461354
$ws.Range("K2").Value = 54.25680821473518

# Row 3 (test_rmse): updated values for all 10 models.
$ws.Range("B3").Value = 0.5793310029159989
$ws.Range("C3").Value = 0.6828392768481389
$ws.Range("D3").Value = 0.7109865587139155
$ws.Range("E3").Value = 0.5587914946739748
$ws.Range("F3").Value = 0.5625940117129977
$ws.Range("G3").Value = 0.6985583934814392
$ws.Range("H3").Value = 0.6817654642258674
$ws.Range("I3").Value = 0.5439349065121836
$ws.Range("J3").Value = 0.6118274768573116
$ws.Range("K3").Value = 0.6374758219755982

# Row 4 (test_score): updated values for all 10 models.
$ws.Range("B4").Value = 0.8425486947018774
$ws.Range("C4").Value = 0.7864791828555581
$ws.Range("D4").Value = 0.8052878858953885
$ws.Range("E4").Value = 0.8661249249735148
$ws.Range("F4").Value = 0.8413550144787181
$ws.Range("G4").Value = 0.7723485920790469
$ws.Range("H4").Value = 0.7580373529393835
$ws.Range("I4").Value = 0.8906154327761845
$ws.Range("J4").Value = 0.7716833961433427
$ws.Range("K4").Value = 0.7486562908830035

### Sheet 2 ###
$ws = $wb.Worksheets.Item(2)

# Header row (row 1): copy formatting from F1 (bordered/bold/centered style)
# into the new header cells G1:K1, then set their values 6..10.
$ws.Range("F1").Copy()
$ws.Range("G1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10

# Row 2 (test_MAPE): updated values for all 10 models.
$ws.Range("B2").Value = 37.36937362360312
$ws.Range("C2").Value = 34.96575798684802
$ws.Range("D2").Value = 38.22814155549965
$ws.Range("E2").Value = 30.77963682037742
$ws.Range("F2").Value = 45.99958864690752
$ws.Range("G2").Value = 37.07694808593751
$ws.Range("H2").Value = 35.39087536173187
$ws.Range("I2").Value = 48.16734599379718
$ws.Range("J2").Value = 37.45436288791731
$ws.Range("K2").Value = 56.78045671557197

# Row 3 (test_rmse): updated values for all 10 models.
$ws.Range("B3").Value = 0.7047525999704891
$ws.Range("C3").Value = 0.649031809840943
$ws.Range("D3").Value = 0.8487868204622602
$ws.Range("E3").Value = 0.5887793081473495
$ws.Range("F3").Value = 0.5812618797255628
$ws.Range("G3").Value = 0.6797799787069773
$ws.Range("H3").Value = 0.622816288230676
$ws.Range("I3").Value = 0.5739376998274894
$ws.Range("J3").Value = 0.6378962721133219
$ws.Range("K3").Value = 0.6561494821936028

# Row 4 (test_score): updated values for all 10 models.
$ws.Range("B4").Value = 0.7669945399419645
$ws.Range("C4").Value = 0.8070986795002415
$ws.Range("D4").Value = 0.7224971509741622
$ws.Range("E4").Value = 0.8513704267667095
$ws.Range("F4").Value = 0.830652097072924
$ws.Range("G4").Value = 0.7844233848021015
$ws.Range("H4").Value = 0.7980712040942353
$ws.Range("I4").Value = 0.8782155894845962
$ws.Range("J4").Value = 0.7518126339517133
$ws.Range("K4").Value = 0.7337153296528807

### Sheet 3 ###
$ws = $wb.Worksheets.Item(3)

# Header row (row 1): copy formatting from F1 (bordered/bold/centered style)
# into the new header cells G1:K1, then set their values 6..10.
$ws.Range("F1").Copy()
$ws.Range("G1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10

# Row 2 (test_MAPE): updated values for all 10 models.
$ws.Range("B2").Value = 37.72818841984554
$ws.Range("C2").Value = 45.97664956604137
$ws.Range("D2").Value = 41.02222074193498
$ws.Range("E2").Value = 35.15516632073809
$ws.Range("F2").Value = 54.8614797959845
$ws.Range("G2").Value = 37.45214207266446
$ws.Range("H2").Value = 39.56353281943321
$ws.Range("I2").Value = 52.25110058982995
$ws.Range("J2").Value = 43.47459823000558
$ws.Range("K2").Value = 58.67119656657275

# Row 3 (test_rmse): updated values for all 10 models.
$ws.Range("B3").Value = 0.6937907726114121
$ws.Range("C3").Value = 0.7378380885492311
$ws.Range("D3").Value = 0.8326363970159318
$ws.Range("E3").Value = 0.6620307310516503
$ws.Range("F3").Value = 0.6371260303711855
$ws.Range("G3").Value = 0.7288098208368927
$ws.Range("H3").Value = 0.6728500317966467
$ws.Range("I3").Value = 0.594097668205498
$ws.Range("J3").Value = 0.6417088637896213
$ws.Range("K3").Value = 0.6467913670257639

# Row 4 (test_score): updated values for all 10 models.
$ws.Range("B4").Value = 0.7741865721242489
$ws.Range("C4").Value = 0.7506982216709186
$ws.Range("D4").Value = 0.732957137976106
$ws.Range("E4").Value = 0.8120871617810239
$ws.Range("F4").Value = 0.7965363473213727
$ws.Range("G4").Value = 0.7522045415465047
$ws.Range("H4").Value = 0.7643242562892085
$ws.Range("I4").Value = 0.8695098025803693
$ws.Range("J4").Value = 0.7488370253075691
$ws.Range("K4").Value = 0.7412567590627132

### Sheet 4 ###
$ws = $wb.Worksheets.Item(4)

# Header row (row 1): copy formatting from F1 (bordered/bold/centered style)
# into the new header cells G1:K1, then set their values 6..10.
$ws.Range("F1").Copy()
$ws.Range("G1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10

# Row 2 (test_MAPE): updated values for all 10 models.
$ws.Range("B2").Value = 49.49359112518707
$ws.Range("C2").Value = 61.96802891958342
$ws.Range("D2").Value = 55.56092945670751
$ws.Range("E2").Value = 54.86142324937349
$ws.Range("F2").Value = 71.22568683739945
$ws.Range("G2").Value = 47.83245738043706
$ws.Range("H2").Value = 54.03652176510341
$ws.Range("I2").Value = 90.69416180347929
$ws.Range("J2").Value = 61.89426987663489
$ws.Range("K2").Value = 72.23383286291525

# Row 3 (test_rmse): updated values for all 10 models.
$ws.Range("B3").Value = 0.8143813899525644
$ws.Range("C3").Value = 0.8648189780504634
$ws.Range("D3").Value = 0.8575473935589115
$ws.Range("E3").Value = 0.8794596946131777
$ws.Range("F3").Value = 0.750950727874952
$ws.Range("G3").Value = 0.8643342777671682
$ws.Range("H3").Value = 0.7963426989774348
$ws.Range("I3").Value = 0.8257256394722737
$ws.Range("J3").Value = 0.7254839807046519
$ws.Range("K3").Value = 0.8365870004528093

# Row 4 (test_score): updated values for all 10 models.
$ws.Range("B4").Value = 0.6888653311640482
$ws.Range("C4").Value = 0.6575054554158579
$ws.Range("D4").Value = 0.7167392153103405
$ws.Range("E4").Value = 0.6683866283527444
$ws.Range("F4").Value = 0.7173434693014356
$ws.Range("G4").Value = 0.651479456005993
$ws.Range("H4").Value = 0.6698750585871875
$ws.Range("I4").Value = 0.7479227034049891
$ws.Range("J4").Value = 0.6789777041598343
$ws.Range("K4").Value = 0.5671246972442818

